$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.278.89'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.282.90'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.80'
$ws.Range("E5").Value = '  +2.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.46'
$ws.Range("E6").Value = '  +0.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.648'
$ws.Range("E7").Value = '  +8.30%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.126'
$ws.Range("E9").Value = '  -2.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.76'
$ws.Range("E10").Value = '  +2.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.405'
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.850.63'
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.131'
$ws.Range("E13").Value = '  -5.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '66.270.04'
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.50'
$ws.Range("E15").Value = '  -1.50%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.364.92'
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000164'
$ws.Range("E17").Value = '  -1.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '435.55'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.25'
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.52'
$ws.Range("E20").Value = '  -2.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.43'
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.14'
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.421.77'
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.510'
$ws.Range("E25").Value = '  -0.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000114'
$ws.Range("E26").Value = '  -3.70%  '
$ws.Range("E27").Value = '  +2.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.89'
$ws.Range("E28").Value = '  -0.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.95'
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.35'
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.20'
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.64'
$ws.Range("E34").Value = '  -1.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.20'
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.90'
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.42'
$ws.Range("E37").Value = '  -4.54%  '
$ws.Range("E38").Value = '  -2.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.28'
$ws.Range("E39").Value = '  -3.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.786.32'
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.774'
$ws.Range("E41").Value = '  -1.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.36'
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.13'
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.28'
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0664'
$ws.Range("E45").Value = '  -1.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.32'
$ws.Range("E46").Value = '  +0.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '321.90'
$ws.Range("E47").Value = '  +1.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.35'
$ws.Range("E48").Value = '  -3.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0269'
$ws.Range("E49").Value = '  -0.46%  '
$ws.Range("E50").Value = '  +4.98%  '
